$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "57.824.78"
$ws.Range("E2").Value = "  +2.59%  "
$ws.Range("D3").Value = "3.064.04"
$ws.Range("E3").Value = "  +2.12%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "'517.74"
$ws.Range("E5").Value = "  +2.11%  "
$ws.Range("D6").Value = "'142.22"
$ws.Range("E6").Value = "  +2.96%  "
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("D9").Value = "'7.30"
$ws.Range("E9").Value = "  +2.40%  "
$ws.Range("E10").Value = "  -0.28%  "
$ws.Range("E11").Value = "  +3.07%  "
$ws.Range("D12").Value = "3.585.99"
$ws.Range("E12").Value = "  +2.21%  "
$ws.Range("E13").Value = "  +3.18%  "
$ws.Range("D14").Value = "'25.97"
$ws.Range("E14").Value = "  +1.41%  "
$ws.Range("E15").Value = "  -0.33%  "
$ws.Range("D16").Value = "57.841.90"
$ws.Range("E16").Value = "  +2.65%  "
$ws.Range("D17").Value = "3.061.40"
$ws.Range("E17").Value = "  +2.23%  "
$ws.Range("D18").Value = "'6.10"
$ws.Range("E18").Value = "  +1.79%  "
$ws.Range("D19").Value = "'12.83"
$ws.Range("E19").Value = "  -0.94%  "
$ws.Range("D20").Value = "'8.12"
$ws.Range("E20").Value = "  +0.56%  "
$ws.Range("D21").Value = "'330.31"
$ws.Range("E21").Value = "  -0.44%  "
$ws.Range("D23").Value = "'0.499"
$ws.Range("E23").Value = "  +0.39%  "
$ws.Range("D24").Value = "'65.66"
$ws.Range("E24").Value = "  +1.05%  "
$ws.Range("E25").Value = "  +2.84%  "
$ws.Range("D27").Value = "0.0₃0903"
$ws.Range("E27").Value = "  -4.52%  "
$ws.Range("E28").Value = "  +0.58%  "
$ws.Range("D29").Value = "'7.21"
$ws.Range("E29").Value = "  +4.22%  "
$ws.Range("E30").Value = "  +2.11%  "
$ws.Range("E31").Value = "  +3.12%  "
$ws.Range("D32").Value = "'20.71"
$ws.Range("E32").Value = "  +2.02%  "
$ws.Range("D33").Value = "'154.74"
$ws.Range("E33").Value = "  +0.65%  "
$ws.Range("D34").Value = "'4.50"
$ws.Range("E34").Value = "  +0.56%  "
$ws.Range("B35").Value = "EnergySwap"
$ws.Range("C35").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D35").Value = "'27.02"
$ws.Range("E35").Value = "  +2.86%  "
$ws.Range("B36").Value = "Aptos"
$ws.Range("C36").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D36").Value = "'5.94"
$ws.Range("E36").Value = "  +1.54%  "
$ws.Range("E37").Value = "  +1.34%  "
$ws.Range("D39").Value = "3.104.61"
$ws.Range("E39").Value = "  +2.27%  "
$ws.Range("E40").Value = "  +3.09%  "
$ws.Range("D41").Value = "'36.57"
$ws.Range("E41").Value = "  -0.40%  "
$ws.Range("E42").Value = "  +0.03%  "
$ws.Range("D43").Value = "'0.655"
$ws.Range("E43").Value = "  +0.57%  "
$ws.Range("D44").Value = "2.257.40"
$ws.Range("E44").Value = "  +3.51%  "
$ws.Range("E45").Value = "  +9.48%  "
$ws.Range("D46").Value = "'20.63"
$ws.Range("E46").Value = "  +5.97%  "
$ws.Range("E47").Value = "  +0.76%  "
$ws.Range("B48").Value = "ONDO"
$ws.Range("C48").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D48").Value = "'0.931"
$ws.Range("E48").Value = "  +0.39%  "
$ws.Range("B49").Value = "Cosmos"
$ws.Range("C49").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D49").Value = "'5.88"
$ws.Range("E49").Value = "  +0.05%  "
$ws.Range("D50").Value = "'0.735"
$ws.Range("E50").Value = "  +8.70%  "
$ws.Range("D51").Value = "'257.92"
$ws.Range("E51").Value = "  +12.08%  "
